$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-09-05 to 2023-09-06
# (serial date 45174 -> 45175), preserving existing cell formatting.
$ws.Range("C2:C5").Value = 45175
